$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalCommitment")

# Fix the CF1/CF2/CF3 headers to have a space: "CF 1", "CF 2", "CF 3"
$ws.Range("S1").Value = "CF 1"
$ws.Range("T1").Value = "CF 2"
$ws.Range("U1").Value = "CF 3"

# Add new investor access columns T (letter code) and U (amount) for rows 2-9
$letters = @("A", "B", "C", "D", "E", "F", "G", "H")
$amounts = @(100, 200, 300, 400, 500, 600, 700, 800)

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 20).Value = $letters[$i]
    $ws.Cells.Item($row, 21).Value = $amounts[$i]
}

# Update the selection to match the edited range
$ws.Range("U2:U9").Select()
